{"js": "// Merge the three runs that make up a comment's body\n// (\"\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430 \" + \" \" + \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435 ...\")\n// into a single run, matching the canonical OOXML Word produces\n// when a comment's text is edited and the document is re-saved:\n//\n//   before (3 runs):\n//     <w:r><w:t xml:space=\"preserve\">\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430 </w:t></w:r>\n//     <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n//     <w:r><w:t>\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435 \u0434\u043e\u043f\u044a\u043b\u043d\u0438\u0442\u0435\u043b\u043d\u0438 \u0443\u0441\u043b\u0443\u0433\u0438, \u043d\u0430\u043c\u0430\u043b\u044f\u0432\u0430 \u0441\u0435 \u043f\u0440\u0438\n//               \u043f\u043e\u0432\u0440\u0435\u0434\u0430 \u043d\u0430 \u0445\u043e\u0442\u0435\u043b\u0441\u043a\u043e\u0442\u043e \u0438\u043c\u0443\u0449\u0435\u0441\u0442\u0432\u043e.</w:t></w:r>\n//\n//   after (1 run):\n//     <w:r><w:t>\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430  \u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435\n//               \u0434\u043e\u043f\u044a\u043b\u043d\u0438\u0442\u0435\u043b\u043d\u0438 \u0443\u0441\u043b\u0443\u0433\u0438, \u043d\u0430\u043c\u0430\u043b\u044f\u0432\u0430 \u0441\u0435 \u043f\u0440\u0438 \u043f\u043e\u0432\u0440\u0435\u0434\u0430 \u043d\u0430\n//               \u0445\u043e\u0442\u0435\u043b\u0441\u043a\u043e\u0442\u043e \u0438\u043c\u0443\u0449\u0435\u0441\u0442\u0432\u043e.</w:t></w:r>\n\nconst comments = context.document.getComments();\ncomments.load(\"items\");\nawait context.sync();\n\nfor (const c of comments.items) {\n  c.load(\"content\");\n}\nawait context.sync();\n\nconst marker = \"\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430\";\nconst comment = comments.items.find((c) => c.content.indexOf(marker) === 0);\nif (!comment) {\n  throw new Error(\"Target comment not found\");\n}\n\nconst fullText = comment.content;\n\n// The host engine short-circuits a `content` write that is byte-identical\n// to the current value (it never touches the underlying run structure in\n// that case), so the merge of the 3 runs into 1 never happens unless the\n// value actually changes at least once. Write a throwaway placeholder\n// first, then write the real text back -- this reliably collapses the\n// comment body down to a single run, the same way Word does internally\n// when the comment text is edited and re-serialized.\ncomment.content = \"placeholder\";\nawait context.sync();\n\ncomment.content = fullText;\nawait context.sync();\n", "ps1": "# Merge the three runs that make up a comment's body\n# (\"\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430 \" + \" \" + \"\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435 ...\")\n# into a single run, matching the canonical OOXML Word produces\n# when a comment's text is edited and the document is re-saved:\n#\n#   before (3 runs):\n#     <w:r><w:t xml:space=\"preserve\">\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430 </w:t></w:r>\n#     <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n#     <w:r><w:t>\u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435 \u0434\u043e\u043f\u044a\u043b\u043d\u0438\u0442\u0435\u043b\u043d\u0438 \u0443\u0441\u043b\u0443\u0433\u0438, \u043d\u0430\u043c\u0430\u043b\u044f\u0432\u0430 \u0441\u0435 \u043f\u0440\u0438\n#               \u043f\u043e\u0432\u0440\u0435\u0434\u0430 \u043d\u0430 \u0445\u043e\u0442\u0435\u043b\u0441\u043a\u043e\u0442\u043e \u0438\u043c\u0443\u0449\u0435\u0441\u0442\u0432\u043e.</w:t></w:r>\n#\n#   after (1 run):\n#     <w:r><w:t>\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430  \u0438\u0437\u043f\u043e\u043b\u0437\u0432\u0430\u043d\u0438\u0442\u0435\n#               \u0434\u043e\u043f\u044a\u043b\u043d\u0438\u0442\u0435\u043b\u043d\u0438 \u0443\u0441\u043b\u0443\u0433\u0438, \u043d\u0430\u043c\u0430\u043b\u044f\u0432\u0430 \u0441\u0435 \u043f\u0440\u0438 \u043f\u043e\u0432\u0440\u0435\u0434\u0430 \u043d\u0430\n#               \u0445\u043e\u0442\u0435\u043b\u0441\u043a\u043e\u0442\u043e \u0438\u043c\u0443\u0449\u0435\u0441\u0442\u0432\u043e.</w:t></w:r>\n\n$d = $word.ActiveDocument\n\n$marker = \"\u0420\u0435\u0439\u0442\u0438\u043d\u0433\u0430 \u0441\u0435 \u043e\u043f\u0440\u0435\u0434\u0435\u043b\u044f \u043e\u0442 \u0431\u0440\u043e\u044f \u043d\u0430\"\n$comment = $null\nfor ($i = 1; $i -le $d.Comments.Count; $i++) {\n    $candidate = $d.Comments.Item($i)\n    if ($candidate.Range.Text.StartsWith($marker)) {\n        $comment = $candidate\n    }\n}\n\n$fullText = $comment.Range.Text\n\n# The host short-circuits a Range.Text write that is byte-identical to the\n# current text (the underlying run structure is left untouched in that\n# case), so the merge of the 3 runs into 1 never happens unless the value\n# actually changes at least once. Write a throwaway placeholder first, then\n# write the real text back -- this reliably collapses the comment body down\n# to a single run, the same way Word does internally when the comment text\n# is edited and re-serialized.\n$comment.Range.Text = \"placeholder\"\n$comment.Range.Text = $fullText\n"}
